$d = $word.ActiveDocument

$replacements = @(
    @{old='922×6=5532'; new='375×2=750'},
    @{old='584×3=1752'; new='162×3=486'},
    @{old='964×9=8676'; new='214×6=1284'},
    @{old='564×9=5076'; new='847×9=7623'},
    @{old='683×2=1366'; new='298×5=1490'},
    @{old='418×9=3762'; new='338×5=1690'},
    @{old='617×7=4319'; new='665×4=2660'},
    @{old='846×2=1692'; new='757×9=6813'},
    @{old='401×8=3208'; new='638×4=2552'},
    @{old='623×9=5607'; new='408×2=816'},
    @{old='924×4=3696'; new='789×3=2367'},
    @{old='406×7=2842'; new='917×6=5502'},
    @{old='123×2=246';  new='883×4=3532'},
    @{old='931×4=3724'; new='533×7=3731'},
    @{old='211×7=1477'; new='479×2=958'},
    @{old='737×6=4422'; new='684×4=2736'},
    @{old='536×9=4824'; new='330×4=1320'},
    @{old='271×2=542';  new='697×4=2788'},
    @{old='353×4=1412'; new='101×9=909'},
    @{old='356×7=2492'; new='899×2=1798'},
    @{old='855×5=4275'; new='982×3=2946'},
    @{old='817×3=2451'; new='258×7=1806'},
    @{old='305×3=915';  new='741×6=4446'},
    @{old='184×4=736';  new='106×7=742'},
    @{old='911×8=7288'; new='853×5=4265'}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
